# Auto-generated edit script applying the cell-value changes described in the diff
# (scheduled market-board profit recalculation for the Balmung_Profits workbook).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 36116384
$ws.Range("I43").Value = 144446540
$ws.Range("J43").Value = 6330.8887
$ws.Range("K43").Value = 144446540
$ws.Range("L43").Value = 6330.8887
$ws.Range("M43").Value = -144446471
$ws.Range("N43").Value = -6468.8887
$ws.Range("H70").Value = 970.8
$ws.Range("I70").Value = 961.6667
$ws.Range("K70").Value = 2885.0001
$ws.Range("M70").Value = -2615.0001
$ws.Range("H73").Value = 970.8
$ws.Range("I73").Value = 961.6667
$ws.Range("K73").Value = 2885.0001
$ws.Range("M73").Value = -1949.0001
$ws.Range("H94").Value = 719
$ws.Range("I94").Value = 719
$ws.Range("K94").Value = 719
$ws.Range("M94").Value = -268
$ws.Range("H100").Value = 3694.5557
$ws.Range("J100").Value = 3995.3333
$ws.Range("L100").Value = 3995.3333
$ws.Range("N100").Value = -5077.3333
$ws.Range("H107").Value = 560.2857
$ws.Range("I107").Value = 486
$ws.Range("K107").Value = 486
$ws.Range("M107").Value = 1434
$ws.Range("H132").Value = 25133.453
$ws.Range("I132").Value = 33414.547
$ws.Range("J132").Value = 1795.8182
$ws.Range("K132").Value = 100243.641
$ws.Range("L132").Value = 5387.4546
$ws.Range("M132").Value = -97713.641
$ws.Range("N132").Value = -10447.4546
$ws.Range("H137").Value = 7694109.5
$ws.Range("I137").Value = 1515.5
$ws.Range("K137").Value = 4546.5
$ws.Range("M137").Value = -1996.5
$ws.Range("H138").Value = 5205.5483
$ws.Range("I138").Value = 7913.885
$ws.Range("K138").Value = 23741.655
$ws.Range("M138").Value = -18601.655
$ws.Range("H141").Value = 2597.1428
$ws.Range("I141").Value = 2146.75
$ws.Range("J141").Value = 5299.5
$ws.Range("K141").Value = 6440.25
$ws.Range("L141").Value = 15898.5
$ws.Range("M141").Value = -1260.25
$ws.Range("N141").Value = -26258.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 1396.3334
$ws.Range("J5").Value = 2000
$ws.Range("L5").Value = 2000
$ws.Range("N5").Value = -2224
$ws.Range("H45").Value = 151665.6
$ws.Range("I45").Value = 151665.6
$ws.Range("K45").Value = 151665.6
$ws.Range("M45").Value = -151288.6
$ws.Range("H61").Value = 1258634.1
$ws.Range("I61").Value = 4811.971
$ws.Range("K61").Value = 4811.971
$ws.Range("M61").Value = -4599.971
$ws.Range("H74").Value = 798093.9399999999
$ws.Range("I74").Value = 1531.4615
$ws.Range("J74").Value = 1488448.1
$ws.Range("K74").Value = 1531.4615
$ws.Range("L74").Value = 1488448.1
$ws.Range("M74").Value = -657.4614999999999
$ws.Range("N74").Value = -1490196.1
$ws.Range("H77").Value = 798093.9399999999
$ws.Range("I77").Value = 1531.4615
$ws.Range("J77").Value = 1488448.1
$ws.Range("K77").Value = 7657.307499999999
$ws.Range("L77").Value = 7442240.5
$ws.Range("M77").Value = -3289.307499999999
$ws.Range("N77").Value = -7450976.5
$ws.Range("H97").Value = 7230.5
$ws.Range("I97").Value = 7645.8667
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 7645.8667
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -7149.8667
$ws.Range("N97").Value = -1992
$ws.Range("H102").Value = 7331.7617
$ws.Range("I102").Value = 7573.35
$ws.Range("K102").Value = 7573.35
$ws.Range("M102").Value = -5951.35
$ws.Range("H132").Value = 1522.4032
$ws.Range("I132").Value = 1335.3448
$ws.Range("K132").Value = 4006.0344
$ws.Range("M132").Value = -1476.0344
$ws.Range("H136").Value = 1258634.1
$ws.Range("I136").Value = 4811.971
$ws.Range("K136").Value = 14435.913
$ws.Range("M136").Value = -11885.913

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 1396.3334
$ws.Range("J4").Value = 2000
$ws.Range("L4").Value = 2000
$ws.Range("N4").Value = -2230
$ws.Range("H20").Value = 1310.3158
$ws.Range("I20").Value = 1114.4
$ws.Range("J20").Value = 1528
$ws.Range("K20").Value = 1114.4
$ws.Range("L20").Value = 1528
$ws.Range("M20").Value = -867.4000000000001
$ws.Range("N20").Value = -2022
$ws.Range("H32").Value = 50000
$ws.Range("J32").Value = 50000
$ws.Range("L32").Value = 50000
$ws.Range("N32").Value = -50768
$ws.Range("H99").Value = 6608.5
$ws.Range("I99").Value = 10265.929
$ws.Range("K99").Value = 10265.929
$ws.Range("M99").Value = -8767.929
$ws.Range("H134").Value = 16982974
$ws.Range("I134").Value = 1726.119
$ws.Range("J134").Value = 81820460
$ws.Range("K134").Value = 5178.357
$ws.Range("L134").Value = 245461380
$ws.Range("M134").Value = -2643.357
$ws.Range("N134").Value = -245466450
$ws.Range("H135").Value = 75998.7
$ws.Range("J135").Value = 75998.7
$ws.Range("L135").Value = 75998.7
$ws.Range("N135").Value = -86138.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1508281.5
$ws.Range("I3").Value = 2996563
$ws.Range("K3").Value = 2996563
$ws.Range("M3").Value = -2996450
$ws.Range("H33").Value = 12000
$ws.Range("I33").Value = 12000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 12000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -11621
$ws.Range("H58").Value = 2404.8965
$ws.Range("I58").Value = 1857.15
$ws.Range("K58").Value = 1857.15
$ws.Range("M58").Value = -1654.15
$ws.Range("H75").Value = 77664.336
$ws.Range("J75").Value = 77664.336
$ws.Range("L75").Value = 77664.336
$ws.Range("N75").Value = -79660.336
$ws.Range("H78").Value = 77664.336
$ws.Range("J78").Value = 77664.336
$ws.Range("L78").Value = 232993.008
$ws.Range("N78").Value = -242977.008
$ws.Range("H97").Value = 69699.75
$ws.Range("J97").Value = 69699.75
$ws.Range("L97").Value = 69699.75
$ws.Range("N97").Value = -71681.75
$ws.Range("H134").Value = 1442.7906
$ws.Range("I134").Value = 1290.4865
$ws.Range("K134").Value = 3871.4595
$ws.Range("M134").Value = -1336.4595
$ws.Range("H136").Value = 2404.8965
$ws.Range("I136").Value = 1857.15
$ws.Range("K136").Value = 5571.450000000001
$ws.Range("M136").Value = -3021.450000000001
$ws.Range("N33").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 4260.8096
$ws.Range("J39").Value = 5320.8125
$ws.Range("L39").Value = 15962.4375
$ws.Range("N39").Value = -16550.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1727.5
$ws.Range("I97").Value = 1455
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 1455
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -959
$ws.Range("N97").Value = -2992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3663.25
$ws.Range("I16").Value = 3417.7778
$ws.Range("J16").Value = 4399.6665
$ws.Range("K16").Value = 3417.7778
$ws.Range("L16").Value = 4399.6665
$ws.Range("M16").Value = -3247.7778
$ws.Range("N16").Value = -4739.6665
$ws.Range("H22").Value = 2963.4119
$ws.Range("I22").Value = 438.9
$ws.Range("K22").Value = 438.9
$ws.Range("M22").Value = -143.9
$ws.Range("H23").Value = 724285.1
$ws.Range("I23").Value = 724285.1
$ws.Range("K23").Value = 724285.1
$ws.Range("M23").Value = -724055.1
$ws.Range("H27").Value = 2963.4119
$ws.Range("I27").Value = 438.9
$ws.Range("K27").Value = 438.9
$ws.Range("M27").Value = -331.9
$ws.Range("H74").Value = 41160.445
$ws.Range("I74").Value = 38634.855
$ws.Range("J74").Value = 50000
$ws.Range("K74").Value = 38634.855
$ws.Range("L74").Value = 50000
$ws.Range("M74").Value = -37636.855
$ws.Range("N74").Value = -51996
$ws.Range("H77").Value = 41160.445
$ws.Range("I77").Value = 38634.855
$ws.Range("J77").Value = 50000
$ws.Range("K77").Value = 115904.565
$ws.Range("L77").Value = 150000
$ws.Range("M77").Value = -110912.565
$ws.Range("N77").Value = -159984
$ws.Range("H99").Value = 46999
$ws.Range("I99").Value = 45665.5
$ws.Range("K99").Value = 45665.5
$ws.Range("M99").Value = -42670.5
$ws.Range("H122").Value = 3824.475
$ws.Range("I122").Value = 3073.926
$ws.Range("J122").Value = 5383.3076
$ws.Range("K122").Value = 9221.778
$ws.Range("L122").Value = 16149.9228
$ws.Range("M122").Value = -6771.778
$ws.Range("N122").Value = -21049.9228
$ws.Range("H132").Value = 4439.8096
$ws.Range("J132").Value = 8096
$ws.Range("L132").Value = 24288
$ws.Range("N132").Value = -29348
$ws.Range("H136").Value = 2378.4363
$ws.Range("I136").Value = 890.55554
$ws.Range("J136").Value = 5197.579
$ws.Range("K136").Value = 2671.66662
$ws.Range("L136").Value = 15592.737
$ws.Range("M136").Value = -121.66662
$ws.Range("N136").Value = -20692.737

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 47775.11
$ws.Range("J75").Value = 47500
$ws.Range("L75").Value = 47500
$ws.Range("N75").Value = -49372
$ws.Range("H78").Value = 47775.11
$ws.Range("J78").Value = 47500
$ws.Range("L78").Value = 142500
$ws.Range("N78").Value = -151860
$ws.Range("H107").Value = 1588266.8
$ws.Range("J107").Value = 7143758
$ws.Range("L107").Value = 21431274
$ws.Range("N107").Value = -21435114
$ws.Range("H132").Value = 15995.869
$ws.Range("I132").Value = 19752.426
$ws.Range("K132").Value = 59257.278
$ws.Range("M132").Value = -56727.278
$ws.Range("H136").Value = 23267.744
$ws.Range("I136").Value = 27519
$ws.Range("K136").Value = 82557
$ws.Range("M136").Value = -80007
